$wb = $excel.ActiveWorkbook

# Sheet1 ("Genetics" data) - row 2, column A: update case id
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2").Value = "CA-II7V08WU"

# Sheet2 (status sheet) - row 3, column B: update status from fail to pass
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B3").Value = "pass"
